$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A135:B135").NumberFormat = "@"
$ws.Range("A135").Value = "134"
$ws.Range("B135").Value = "2020/07/16"
$ws.Range("C135").Value = 9546.0
$ws.Range("D135").Value = 578.74
$ws.Range("E135").Value = 582.86

$ws.Range("A136:B136").NumberFormat = "@"
$ws.Range("A136").Value = "135"
$ws.Range("B136").Value = "2020/07/17"
$ws.Range("C136").Value = 9969.0
$ws.Range("D136").Value = 578.91
$ws.Range("E136").Value = 584.16

$ws.Range("A137:B137").NumberFormat = "@"
$ws.Range("A137").Value = "136"
$ws.Range("B137").Value = "2020/07/18"
$ws.Range("C137").Value = 10551.0
$ws.Range("D137").Value = 578.7
$ws.Range("E137").Value = 584.52

$ws.Range("A138:B138").NumberFormat = "@"
$ws.Range("A138").Value = "137"
$ws.Range("B138").Value = "2020/07/19"
$ws.Range("C138").Value = 11114.0
$ws.Range("D138").Value = 578.7
$ws.Range("E138").Value = 584.52

$ws.Range("A139:B139").NumberFormat = "@"
$ws.Range("A139").Value = "138"
$ws.Range("B139").Value = "2020/07/20"
$ws.Range("C139").Value = 11534.0
$ws.Range("D139").Value = 578.7
$ws.Range("E139").Value = 584.52

$ws.Range("A140:B140").NumberFormat = "@"
$ws.Range("A140").Value = "139"
$ws.Range("B140").Value = "2020/07/21"
$ws.Range("C140").Value = 11811.0
$ws.Range("D140").Value = 576.47
$ws.Range("E140").Value = 584.94

$ws.Range("A141:B141").NumberFormat = "@"
$ws.Range("A141").Value = "140"
$ws.Range("B141").Value = "2020/07/22"
$ws.Range("C141").Value = 12361.0
$ws.Range("D141").Value = 576.73
$ws.Range("E141").Value = 585.21

$ws.Range("A142:B142").NumberFormat = "@"
$ws.Range("A142").Value = "141"
$ws.Range("B142").Value = "2020/07/23"
$ws.Range("C142").Value = 13129.0
$ws.Range("D142").Value = 579.78
$ws.Range("E142").Value = 585.62

$ws.Range("A143:B143").NumberFormat = "@"
$ws.Range("A143").Value = "142"
$ws.Range("B143").Value = "2020/07/24"
$ws.Range("C143").Value = 13669.0
$ws.Range("D143").Value = 579.6
$ws.Range("E143").Value = 586.01

$ws.Range("A144:B144").NumberFormat = "@"
$ws.Range("A144").Value = "143"
$ws.Range("B144").Value = "2020/07/25"
$ws.Range("C144").Value = 14600.0
$ws.Range("D144").Value = 579.69
$ws.Range("E144").Value = 586.27

$ws.Range("A145:B145").NumberFormat = "@"
$ws.Range("A145").Value = "144"
$ws.Range("B145").Value = "2020/07/26"
$ws.Range("C145").Value = 15229.0
$ws.Range("D145").Value = 579.69
$ws.Range("E145").Value = 586.27

$ws.Range("A146:B146").NumberFormat = "@"
$ws.Range("A146").Value = "145"
$ws.Range("B146").Value = "2020/07/27"
$ws.Range("C146").Value = 15841.0
$ws.Range("D146").Value = 579.69
$ws.Range("E146").Value = 586.27

$ws.Range("A147:B147").NumberFormat = "@"
$ws.Range("A147").Value = "146"
$ws.Range("B147").Value = "2020/07/28"
$ws.Range("C147").Value = 16344.0
$ws.Range("D147").Value = 579.69
$ws.Range("E147").Value = 586.27

$ws.Range("A148:B148").NumberFormat = "@"
$ws.Range("A148").Value = "147"
$ws.Range("B148").Value = "2020/07/29"
$ws.Range("C148").Value = 16800.0
$ws.Range("D148").Value = 580.49
$ws.Range("E148").Value = 586.9

$ws.Range("A149:B149").NumberFormat = "@"
$ws.Range("A149").Value = "148"
$ws.Range("B149").Value = "2020/07/30"
$ws.Range("C149").Value = 17290.0
$ws.Range("D149").Value = 580.51
$ws.Range("E149").Value = 586.91

$ws.Range("A150:B150").NumberFormat = "@"
$ws.Range("A150").Value = "149"
$ws.Range("B150").Value = "2020/07/31"
$ws.Range("C150").Value = 17820.0
$ws.Range("D150").Value = 582.15
$ws.Range("E150").Value = 588.33

$ws.Range("A151:B151").NumberFormat = "@"
$ws.Range("A151").Value = "150"
$ws.Range("B151").Value = "2020/08/01"
$ws.Range("C151").Value = 18187.0
$ws.Range("D151").Value = 582.41
$ws.Range("E151").Value = 590.74

$ws.Range("A152:B152").NumberFormat = "@"
$ws.Range("A152").Value = "151"
$ws.Range("B152").Value = "2020/08/02"
$ws.Range("C152").Value = 18975.0
$ws.Range("D152").Value = 582.41
$ws.Range("E152").Value = 590.74

$ws.Range("A153:B153").NumberFormat = "@"
$ws.Range("A153").Value = "152"
$ws.Range("B153").Value = "2020/08/03"
$ws.Range("C153").Value = 19402.0
$ws.Range("D153").Value = 582.41
$ws.Range("E153").Value = 590.74

$ws.Range("A154:B154").NumberFormat = "@"
$ws.Range("A154").Value = "153"
$ws.Range("B154").Value = "2020/08/04"
$ws.Range("C154").Value = 19837.0
$ws.Range("D154").Value = 584.23
$ws.Range("E154").Value = 592.22

$ws.Range("A155:B155").NumberFormat = "@"
$ws.Range("A155").Value = "154"
$ws.Range("B155").Value = "2020/08/05"
$ws.Range("C155").Value = 20417.0
$ws.Range("D155").Value = 586.65
$ws.Range("E155").Value = 594.44

$ws.Range("A156:B156").NumberFormat = "@"
$ws.Range("A156").Value = "155"
$ws.Range("B156").Value = "2020/08/06"
$ws.Range("C156").Value = 21070.0
$ws.Range("D156").Value = 588.12
$ws.Range("E156").Value = 595.21

$ws.Range("A157:B157").NumberFormat = "@"
$ws.Range("A157").Value = "156"
$ws.Range("B157").Value = "2020/08/07"
$ws.Range("C157").Value = 22081.0
$ws.Range("D157").Value = 588.7
$ws.Range("E157").Value = 596.17

$ws.Range("A158:B158").NumberFormat = "@"
$ws.Range("A158").Value = "157"
$ws.Range("B158").Value = "2020/08/08"
$ws.Range("C158").Value = 22802.0
$ws.Range("D158").Value = 588.22
$ws.Range("E158").Value = 595.91

$ws.Range("A159:B159").NumberFormat = "@"
$ws.Range("A159").Value = "158"
$ws.Range("B159").Value = "2020/08/09"
$ws.Range("C159").Value = 23286.0
$ws.Range("D159").Value = 588.22
$ws.Range("E159").Value = 595.91

$ws.Range("A160:B160").NumberFormat = "@"
$ws.Range("A160").Value = "159"
$ws.Range("B160").Value = "2020/08/10"
$ws.Range("C160").Value = 23872.0
$ws.Range("D160").Value = 588.22
$ws.Range("E160").Value = 595.91

$ws.Range("A161:B161").NumberFormat = "@"
$ws.Range("A161").Value = "160"
$ws.Range("B161").Value = "2020/08/11"
$ws.Range("C161").Value = 24508.0
$ws.Range("D161").Value = 590.69
$ws.Range("E161").Value = 598.21

$ws.Range("A162:B162").NumberFormat = "@"
$ws.Range("A162").Value = "161"
$ws.Range("B162").Value = "2020/08/12"
$ws.Range("C162").Value = 25057.0
$ws.Range("D162").Value = 592.54
$ws.Range("E162").Value = 599.13
